# Update "想去人数" (interested-count) figures for three events that are
# listed on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 160
$wsExhibit.Range("F6").Value = 9347
$wsExhibit.Range("F10").Value = 1117

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 160
$wsAll.Range("F7").Value = 9347
$wsAll.Range("F11").Value = 1117
